$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28-133 down to 29-134.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new record's data.
$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "Macroferia Regional de Talca"
$ws.Range("C28").Value = "Maule"
$ws.Range("D28").Value = 44459
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = 100112045
$ws.Range("G28").Value = "Zapallo"
$ws.Range("H28").Value = "Paine"
$ws.Range("I28").Value = "1a (guarda)"
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 150
$ws.Range("L28").Value = 150
$ws.Range("M28").Value = 150
$ws.Range("N28").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O28").Value = "Región del Maule"
$ws.Range("P28").Value = 150
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = "Hortaliza"
